# Update cryptocurrency price/volume data per the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (the sheet keeps these price/volume columns as text, not numbers),
    # then reset the style so no stray quote-prefix/number-format leaks in.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "37.827.73"
Set-TextCell $ws.Range("E2") "  +1.54%  "
Set-TextCell $ws.Range("D3") "2.048.61"
Set-TextCell $ws.Range("E3") "  +0.99%  "
Set-TextCell $ws.Range("E4") "  +0.09%  "
Set-TextCell $ws.Range("D5") "229.52"
Set-TextCell $ws.Range("E5") "  +0.89%  "
Set-TextCell $ws.Range("D6") "0.616"
Set-TextCell $ws.Range("E6") "  +0.94%  "
Set-TextCell $ws.Range("D7") "58.20"
Set-TextCell $ws.Range("E7") "  +5.30%  "
Set-TextCell $ws.Range("E8") "  -0.01%  "
Set-TextCell $ws.Range("E9") "  +2.08%  "
Set-TextCell $ws.Range("E10") "  +2.14%  "
Set-TextCell $ws.Range("E11") "  +0.62%  "
Set-TextCell $ws.Range("D12") "2.348.89"
Set-TextCell $ws.Range("E12") "  +0.87%  "
Set-TextCell $ws.Range("D13") "14.56"
Set-TextCell $ws.Range("E13") "  +2.09%  "
Set-TextCell $ws.Range("D14") "20.68"
Set-TextCell $ws.Range("E14") "  +1.31%  "
Set-TextCell $ws.Range("E15") "  +1.83%  "
Set-TextCell $ws.Range("E16") "  +0.39%  "
Set-TextCell $ws.Range("D17") "2.047.70"
Set-TextCell $ws.Range("E17") "  +1.39%  "
Set-TextCell $ws.Range("D18") "37.784.65"
Set-TextCell $ws.Range("E18") "  +1.52%  "
Set-TextCell $ws.Range("D19") "6.15"
Set-TextCell $ws.Range("E19") "  -2.17%  "
Set-TextCell $ws.Range("D20") "69.64"
Set-TextCell $ws.Range("E20") "  +0.55%  "
Set-TextCell $ws.Range("D21") "0.0₃0831"
Set-TextCell $ws.Range("E21") "  +1.13%  "
Set-TextCell $ws.Range("D22") "224.35"
Set-TextCell $ws.Range("E22") "  -0.15%  "
Set-TextCell $ws.Range("E23") "  -0.02%  "
Set-TextCell $ws.Range("E24") "  +0.25%  "
Set-TextCell $ws.Range("E25") "  +1.65%  "
Set-TextCell $ws.Range("D26") "166.73"
Set-TextCell $ws.Range("E26") "  +0.82%  "
Set-TextCell $ws.Range("D27") "9.26"
Set-TextCell $ws.Range("E27") "  -0.84%  "
Set-TextCell $ws.Range("E28") "  +3.19%  "
Set-TextCell $ws.Range("E29") "  +0.90%  "
Set-TextCell $ws.Range("E30") "  -0.65%  "
Set-TextCell $ws.Range("E31") "  +1.54%  "
Set-TextCell $ws.Range("E32") "  -0.13%  "
Set-TextCell $ws.Range("E33") "  +13.62%  "
Set-TextCell $ws.Range("D34") "4.58"
Set-TextCell $ws.Range("E34") "  +2.47%  "
Set-TextCell $ws.Range("D35") "0.0612"
Set-TextCell $ws.Range("E35") "  -1.03%  "
Set-TextCell $ws.Range("D36") "2.34"
Set-TextCell $ws.Range("E36") "  -1.35%  "
Set-TextCell $ws.Range("D37") "5.98"
Set-TextCell $ws.Range("E37") "  +9.43%  "
Set-TextCell $ws.Range("E38") "  +4.41%  "
Set-TextCell $ws.Range("E39") "  -0.10%  "
Set-TextCell $ws.Range("E40") "  +0.30%  "
Set-TextCell $ws.Range("D41") "1.484.75"
Set-TextCell $ws.Range("E41") "  +0.39%  "
Set-TextCell $ws.Range("D42") "97.42"
Set-TextCell $ws.Range("E42") "  +1.44%  "
Set-TextCell $ws.Range("E43") "  +2.74%  "
Set-TextCell $ws.Range("D44") "0.0934"
Set-TextCell $ws.Range("E44") "  +1.20%  "
Set-TextCell $ws.Range("D45") "16.57"
Set-TextCell $ws.Range("E45") "  +0.01%  "
Set-TextCell $ws.Range("D46") "4.20"
Set-TextCell $ws.Range("E46") "  +16.62%  "
Set-TextCell $ws.Range("E47") "  -0.53%  "
Set-TextCell $ws.Range("E48") "  -0.55%  "
Set-TextCell $ws.Range("D49") "2.95"
Set-TextCell $ws.Range("E49") "  +0.91%  "
Set-TextCell $ws.Range("D50") "6.99"
Set-TextCell $ws.Range("E50") "  -3.70%  "
Set-TextCell $ws.Range("D51") "2.241.77"
Set-TextCell $ws.Range("E51") "  +1.23%  "
